# TC for Product search added
# Updates the RTM_ALL sheet's Cart (FR_CART_*) rows:
#  - rewrites the requirement-description text for FR_CART_01..FR_CART_05
#    to the new, clearer wording
#  - inserts a new row for a 6th cart requirement (FR_CART_06 / TS_CART_06)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RTM_ALL")
$ws.Activate()

# --- Reword the existing Cart requirement descriptions (column B) ---
$ws.Cells.Item(31, 2).Value = "User should be able to add product to cart"
$ws.Cells.Item(32, 2).Value = "User should be able to add product to cart"
$ws.Cells.Item(33, 2).Value = "User should be able to update product quantity"
$ws.Cells.Item(34, 2).Value = "User should be able to update product quantity"
$ws.Cells.Item(35, 2).Value = "Cart should show correct total price"

# --- Insert a new row 36 for FR_CART_06 / TS_CART_06 (shifts old rows 36+ down by one) ---
$ws.Rows.Item(36).Insert()

$ws.Cells.Item(36, 1).Value = "FR_CART_06"
$ws.Cells.Item(36, 2).Value = "Cart should show correct total price"
$ws.Cells.Item(36, 3).Value = "Add to Cart"
$ws.Cells.Item(36, 4).Value = "TS_CART_06"

# --- Match the saved view state: scrolled down a bit, D35:D36 selected ---
$ws.Range("D35:D36").Select() | Out-Null
